$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.148.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.241.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.77%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.10%  "
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.40"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.28%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.19%  "
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.577.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.868"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.242.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.118.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0988"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.51%  "
$ws.Range("E29").Value = "  -5.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.80%  "
$ws.Range("E32").Value = "  -6.92%  "
$ws.Range("E33").Value = "  -7.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0720"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.62%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0268"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.70%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.100"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.190"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.55%  "
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  +7.71%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.50%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.32%  "
$ws.Range("E51").Value = "  -6.21%  "
